# ============================================================
# edit.ps1
# Implements the commit: "feat: add 2022-Q3 data"
#
# Before: 总计 (sheet1), 2022-Q2 (sheet2), 2022-Q1 (sheet3)
# After:  总计 (sheet1), 2022-Q3 (NEW sheet2), 2022-Q2 (sheet3), 2022-Q1 (sheet4)
#
# 1) Insert a brand-new worksheet named '2022-Q3' right before the
#    existing '2022-Q2' sheet, and fill it with the Q3 fund-holding data.
# 2) Insert a new row into the '总计' (totals) summary sheet for the
#    2022-Q3 totals (16 holdings, 2.75 billion), above the existing
#    2022-Q2 / 2022-Q1 rows.
# ============================================================

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------
# Step 1: insert the new '2022-Q3' worksheet before '2022-Q2'
# ------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)        # "总计"
$wsQ2Old = $wb.Worksheets.Item(2)        # currently "2022-Q2" (will stay "2022-Q2")

$wsQ3 = $wb.Worksheets.Add($wsQ2Old)
$wsQ3.Name = "2022-Q3"

# -- header row (row 1), columns B..H, bold/centered/boxed like the
#    other quarter sheets -- copy the formatting from an existing header cell
$wsQ2Old.Cells.Item(1, 2).Copy()
$headerRange = $wsQ3.Range("B1:H1")
$headerRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$headers = @('基金代码', '基金名称', '基金规模', '股票总仓位', '仓位占比', '持有市值(亿元)', '仓位排名')
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsQ3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# -- data rows 2..17
# columns: A=index(n), B=code(text), C=name(text), D=size(text),
#          E=position(text), F=ratio(text), G=marketvalue(text, except
#          the very last row which is a genuine 0 number), H=rank(n)

# row 2: fund 110023
$wsQ3.Cells.Item(2, 1).Value = 0
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(2, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(2, 1).Value = 0
$wsQ3.Cells.Item(2, 2).NumberFormat = "@"
$wsQ3.Cells.Item(2, 2).Value = '110023'
$wsQ3.Cells.Item(2, 3).Value = '易方达医疗保健行业混合'
$wsQ3.Cells.Item(2, 4).NumberFormat = "@"
$wsQ3.Cells.Item(2, 4).Value = '33.04'
$wsQ3.Cells.Item(2, 5).NumberFormat = "@"
$wsQ3.Cells.Item(2, 5).Value = '90.58'
$wsQ3.Cells.Item(2, 6).NumberFormat = "@"
$wsQ3.Cells.Item(2, 6).Value = '3.08'
$wsQ3.Cells.Item(2, 7).NumberFormat = "@"
$wsQ3.Cells.Item(2, 7).Value = '1.0176'
$wsQ3.Cells.Item(2, 8).Value = 9

# row 3: fund 590002
$wsQ3.Cells.Item(3, 1).Value = 1
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(3, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(3, 1).Value = 1
$wsQ3.Cells.Item(3, 2).NumberFormat = "@"
$wsQ3.Cells.Item(3, 2).Value = '590002'
$wsQ3.Cells.Item(3, 3).Value = '中邮核心成长混合'
$wsQ3.Cells.Item(3, 4).NumberFormat = "@"
$wsQ3.Cells.Item(3, 4).Value = '37.60'
$wsQ3.Cells.Item(3, 5).NumberFormat = "@"
$wsQ3.Cells.Item(3, 5).Value = '73.26'
$wsQ3.Cells.Item(3, 6).NumberFormat = "@"
$wsQ3.Cells.Item(3, 6).Value = '2.10'
$wsQ3.Cells.Item(3, 7).NumberFormat = "@"
$wsQ3.Cells.Item(3, 7).Value = '0.7896'
$wsQ3.Cells.Item(3, 8).Value = 10

# row 4: fund 160921
$wsQ3.Cells.Item(4, 1).Value = 2
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(4, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(4, 1).Value = 2
$wsQ3.Cells.Item(4, 2).NumberFormat = "@"
$wsQ3.Cells.Item(4, 2).Value = '160921'
$wsQ3.Cells.Item(4, 3).Value = '大成多策略混合（LOF）A'
$wsQ3.Cells.Item(4, 4).NumberFormat = "@"
$wsQ3.Cells.Item(4, 4).Value = '2.90'
$wsQ3.Cells.Item(4, 5).NumberFormat = "@"
$wsQ3.Cells.Item(4, 5).Value = '87.66'
$wsQ3.Cells.Item(4, 6).NumberFormat = "@"
$wsQ3.Cells.Item(4, 6).Value = '6.13'
$wsQ3.Cells.Item(4, 7).NumberFormat = "@"
$wsQ3.Cells.Item(4, 7).Value = '0.1778'
$wsQ3.Cells.Item(4, 8).Value = 3

# row 5: fund 090020
$wsQ3.Cells.Item(5, 1).Value = 3
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(5, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(5, 1).Value = 3
$wsQ3.Cells.Item(5, 2).NumberFormat = "@"
$wsQ3.Cells.Item(5, 2).Value = '090020'
$wsQ3.Cells.Item(5, 3).Value = '大成健康产业混合A'
$wsQ3.Cells.Item(5, 4).NumberFormat = "@"
$wsQ3.Cells.Item(5, 4).Value = '2.38'
$wsQ3.Cells.Item(5, 5).NumberFormat = "@"
$wsQ3.Cells.Item(5, 5).Value = '92.81'
$wsQ3.Cells.Item(5, 6).NumberFormat = "@"
$wsQ3.Cells.Item(5, 6).Value = '6.31'
$wsQ3.Cells.Item(5, 7).NumberFormat = "@"
$wsQ3.Cells.Item(5, 7).Value = '0.1502'
$wsQ3.Cells.Item(5, 8).Value = 5

# row 6: fund 090016
$wsQ3.Cells.Item(6, 1).Value = 4
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(6, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(6, 1).Value = 4
$wsQ3.Cells.Item(6, 2).NumberFormat = "@"
$wsQ3.Cells.Item(6, 2).Value = '090016'
$wsQ3.Cells.Item(6, 3).Value = '大成消费主题混合'
$wsQ3.Cells.Item(6, 4).NumberFormat = "@"
$wsQ3.Cells.Item(6, 4).Value = '3.59'
$wsQ3.Cells.Item(6, 5).NumberFormat = "@"
$wsQ3.Cells.Item(6, 5).Value = '90.85'
$wsQ3.Cells.Item(6, 6).NumberFormat = "@"
$wsQ3.Cells.Item(6, 6).Value = '3.72'
$wsQ3.Cells.Item(6, 7).NumberFormat = "@"
$wsQ3.Cells.Item(6, 7).Value = '0.1335'
$wsQ3.Cells.Item(6, 8).Value = 9

# row 7: fund 012045
$wsQ3.Cells.Item(7, 1).Value = 5
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(7, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(7, 1).Value = 5
$wsQ3.Cells.Item(7, 2).NumberFormat = "@"
$wsQ3.Cells.Item(7, 2).Value = '012045'
$wsQ3.Cells.Item(7, 3).Value = '大成医药健康股票A'
$wsQ3.Cells.Item(7, 4).NumberFormat = "@"
$wsQ3.Cells.Item(7, 4).Value = '1.95'
$wsQ3.Cells.Item(7, 5).NumberFormat = "@"
$wsQ3.Cells.Item(7, 5).Value = '92.51'
$wsQ3.Cells.Item(7, 6).NumberFormat = "@"
$wsQ3.Cells.Item(7, 6).Value = '6.30'
$wsQ3.Cells.Item(7, 7).NumberFormat = "@"
$wsQ3.Cells.Item(7, 7).Value = '0.1228'
$wsQ3.Cells.Item(7, 8).Value = 5

# row 8: fund 001898
$wsQ3.Cells.Item(8, 1).Value = 6
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(8, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(8, 1).Value = 6
$wsQ3.Cells.Item(8, 2).NumberFormat = "@"
$wsQ3.Cells.Item(8, 2).Value = '001898'
$wsQ3.Cells.Item(8, 3).Value = '易方达大健康主题灵活配置混合'
$wsQ3.Cells.Item(8, 4).NumberFormat = "@"
$wsQ3.Cells.Item(8, 4).Value = '3.77'
$wsQ3.Cells.Item(8, 5).NumberFormat = "@"
$wsQ3.Cells.Item(8, 5).Value = '89.93'
$wsQ3.Cells.Item(8, 6).NumberFormat = "@"
$wsQ3.Cells.Item(8, 6).Value = '3.19'
$wsQ3.Cells.Item(8, 7).NumberFormat = "@"
$wsQ3.Cells.Item(8, 7).Value = '0.1203'
$wsQ3.Cells.Item(8, 8).Value = 7

# row 9: fund 016062
$wsQ3.Cells.Item(9, 1).Value = 7
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(9, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(9, 1).Value = 7
$wsQ3.Cells.Item(9, 2).NumberFormat = "@"
$wsQ3.Cells.Item(9, 2).Value = '016062'
$wsQ3.Cells.Item(9, 3).Value = '大成多策略混合（LOF）C'
$wsQ3.Cells.Item(9, 4).NumberFormat = "@"
$wsQ3.Cells.Item(9, 4).Value = '1.37'
$wsQ3.Cells.Item(9, 5).NumberFormat = "@"
$wsQ3.Cells.Item(9, 5).Value = '87.66'
$wsQ3.Cells.Item(9, 6).NumberFormat = "@"
$wsQ3.Cells.Item(9, 6).Value = '6.13'
$wsQ3.Cells.Item(9, 7).NumberFormat = "@"
$wsQ3.Cells.Item(9, 7).Value = '0.0840'
$wsQ3.Cells.Item(9, 8).Value = 3

# row 10: fund 002319
$wsQ3.Cells.Item(10, 1).Value = 8
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(10, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(10, 1).Value = 8
$wsQ3.Cells.Item(10, 2).NumberFormat = "@"
$wsQ3.Cells.Item(10, 2).Value = '002319'
$wsQ3.Cells.Item(10, 3).Value = '大成一带一路灵活配置混合'
$wsQ3.Cells.Item(10, 4).NumberFormat = "@"
$wsQ3.Cells.Item(10, 4).Value = '1.26'
$wsQ3.Cells.Item(10, 5).NumberFormat = "@"
$wsQ3.Cells.Item(10, 5).Value = '89.65'
$wsQ3.Cells.Item(10, 6).NumberFormat = "@"
$wsQ3.Cells.Item(10, 6).Value = '4.49'
$wsQ3.Cells.Item(10, 7).NumberFormat = "@"
$wsQ3.Cells.Item(10, 7).Value = '0.0566'
$wsQ3.Cells.Item(10, 8).Value = 8

# row 11: fund 014121
$wsQ3.Cells.Item(11, 1).Value = 9
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(11, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(11, 1).Value = 9
$wsQ3.Cells.Item(11, 2).NumberFormat = "@"
$wsQ3.Cells.Item(11, 2).Value = '014121'
$wsQ3.Cells.Item(11, 3).Value = '大成品质医疗股票A'
$wsQ3.Cells.Item(11, 4).NumberFormat = "@"
$wsQ3.Cells.Item(11, 4).Value = '0.57'
$wsQ3.Cells.Item(11, 5).NumberFormat = "@"
$wsQ3.Cells.Item(11, 5).Value = '91.58'
$wsQ3.Cells.Item(11, 6).NumberFormat = "@"
$wsQ3.Cells.Item(11, 6).Value = '6.96'
$wsQ3.Cells.Item(11, 7).NumberFormat = "@"
$wsQ3.Cells.Item(11, 7).Value = '0.0397'
$wsQ3.Cells.Item(11, 8).Value = 3

# row 12: fund 001365
$wsQ3.Cells.Item(12, 1).Value = 10
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(12, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(12, 1).Value = 10
$wsQ3.Cells.Item(12, 2).NumberFormat = "@"
$wsQ3.Cells.Item(12, 2).Value = '001365'
$wsQ3.Cells.Item(12, 3).Value = '大成正向回报灵活配置混合'
$wsQ3.Cells.Item(12, 4).NumberFormat = "@"
$wsQ3.Cells.Item(12, 4).Value = '0.57'
$wsQ3.Cells.Item(12, 5).NumberFormat = "@"
$wsQ3.Cells.Item(12, 5).Value = '92.44'
$wsQ3.Cells.Item(12, 6).NumberFormat = "@"
$wsQ3.Cells.Item(12, 6).Value = '6.70'
$wsQ3.Cells.Item(12, 7).NumberFormat = "@"
$wsQ3.Cells.Item(12, 7).Value = '0.0382'
$wsQ3.Cells.Item(12, 8).Value = 4

# row 13: fund 012046
$wsQ3.Cells.Item(13, 1).Value = 11
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(13, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(13, 1).Value = 11
$wsQ3.Cells.Item(13, 2).NumberFormat = "@"
$wsQ3.Cells.Item(13, 2).Value = '012046'
$wsQ3.Cells.Item(13, 3).Value = '大成医药健康股票C'
$wsQ3.Cells.Item(13, 4).NumberFormat = "@"
$wsQ3.Cells.Item(13, 4).Value = '0.22'
$wsQ3.Cells.Item(13, 5).NumberFormat = "@"
$wsQ3.Cells.Item(13, 5).Value = '92.51'
$wsQ3.Cells.Item(13, 6).NumberFormat = "@"
$wsQ3.Cells.Item(13, 6).Value = '6.30'
$wsQ3.Cells.Item(13, 7).NumberFormat = "@"
$wsQ3.Cells.Item(13, 7).Value = '0.0139'
$wsQ3.Cells.Item(13, 8).Value = 5

# row 14: fund 014122
$wsQ3.Cells.Item(14, 1).Value = 12
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(14, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(14, 1).Value = 12
$wsQ3.Cells.Item(14, 2).NumberFormat = "@"
$wsQ3.Cells.Item(14, 2).Value = '014122'
$wsQ3.Cells.Item(14, 3).Value = '大成品质医疗股票C'
$wsQ3.Cells.Item(14, 4).NumberFormat = "@"
$wsQ3.Cells.Item(14, 4).Value = '0.08'
$wsQ3.Cells.Item(14, 5).NumberFormat = "@"
$wsQ3.Cells.Item(14, 5).Value = '91.58'
$wsQ3.Cells.Item(14, 6).NumberFormat = "@"
$wsQ3.Cells.Item(14, 6).Value = '6.96'
$wsQ3.Cells.Item(14, 7).NumberFormat = "@"
$wsQ3.Cells.Item(14, 7).Value = '0.0056'
$wsQ3.Cells.Item(14, 8).Value = 3

# row 15: fund 015655
$wsQ3.Cells.Item(15, 1).Value = 13
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(15, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(15, 1).Value = 13
$wsQ3.Cells.Item(15, 2).NumberFormat = "@"
$wsQ3.Cells.Item(15, 2).Value = '015655'
$wsQ3.Cells.Item(15, 3).Value = '富荣医药健康混合A'
$wsQ3.Cells.Item(15, 4).NumberFormat = "@"
$wsQ3.Cells.Item(15, 4).Value = '0.13'
$wsQ3.Cells.Item(15, 5).NumberFormat = "@"
$wsQ3.Cells.Item(15, 5).Value = '82.11'
$wsQ3.Cells.Item(15, 6).NumberFormat = "@"
$wsQ3.Cells.Item(15, 6).Value = '1.61'
$wsQ3.Cells.Item(15, 7).NumberFormat = "@"
$wsQ3.Cells.Item(15, 7).Value = '0.0021'
$wsQ3.Cells.Item(15, 8).Value = 10

# row 16: fund 015656
$wsQ3.Cells.Item(16, 1).Value = 14
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(16, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(16, 1).Value = 14
$wsQ3.Cells.Item(16, 2).NumberFormat = "@"
$wsQ3.Cells.Item(16, 2).Value = '015656'
$wsQ3.Cells.Item(16, 3).Value = '富荣医药健康混合C'
$wsQ3.Cells.Item(16, 4).NumberFormat = "@"
$wsQ3.Cells.Item(16, 4).Value = '0.01'
$wsQ3.Cells.Item(16, 5).NumberFormat = "@"
$wsQ3.Cells.Item(16, 5).Value = '82.11'
$wsQ3.Cells.Item(16, 6).NumberFormat = "@"
$wsQ3.Cells.Item(16, 6).Value = '1.61'
$wsQ3.Cells.Item(16, 7).NumberFormat = "@"
$wsQ3.Cells.Item(16, 7).Value = '0.0002'
$wsQ3.Cells.Item(16, 8).Value = 10

# row 17: fund 016060
$wsQ3.Cells.Item(17, 1).Value = 15
$wsQ2Old.Cells.Item(2, 1).Copy()
$wsQ3.Cells.Item(17, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsQ3.Cells.Item(17, 1).Value = 15
$wsQ3.Cells.Item(17, 2).NumberFormat = "@"
$wsQ3.Cells.Item(17, 2).Value = '016060'
$wsQ3.Cells.Item(17, 3).Value = '大成健康产业混合C'
$wsQ3.Cells.Item(17, 4).NumberFormat = "@"
$wsQ3.Cells.Item(17, 4).Value = '0.00'
$wsQ3.Cells.Item(17, 5).NumberFormat = "@"
$wsQ3.Cells.Item(17, 5).Value = '92.81'
$wsQ3.Cells.Item(17, 6).NumberFormat = "@"
$wsQ3.Cells.Item(17, 6).Value = '6.31'
$wsQ3.Cells.Item(17, 7).Value = 0
$wsQ3.Cells.Item(17, 8).Value = 5

# ------------------------------------------------------------
# Step 2: insert the 2022-Q3 summary row in the '总计' sheet
# ------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Cells.Item(4, 1).Copy()
$wsTotal.Cells.Item(2, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 16
$wsTotal.Cells.Item(2, 4).Value = 2.75

# bump the A-index of the (now shifted down) old rows back to 0,1 sequence
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(4, 1).Value = 2

# ------------------------------------------------------------
# Restore the originally-selected sheet (last sheet, '2022-Q1')
# as the active tab, matching the source workbook's saved view.
# ------------------------------------------------------------
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLast.Activate()

Write-Host "2022-Q3 sheet inserted and 总计 updated."
